# 15 API Access Control - Effective.xlsx
#
# 1. Rename "Sheet1" -> "Data" (Excel automatically keeps the
#    _xlnm._FilterDatabase defined name in sync, updating
#    "Sheet1!$A$1:$J$1001" to "Data!$A$1:$J$1001").
# 2. On the "Legend" sheet, insert a new header row above the existing
#    0..5 / effectiveness-label rows, label the two header cells
#    "Column1" / "Column2", and promote the A1:B7 range to a native
#    Excel Table ("Table1") with those two columns.

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("Sheet1")
$dataSheet.Name = "Data"

$legend = $wb.Worksheets.Item("Legend")
$legend.Activate()

# Push the existing data down one row to make room for the table header.
$legend.Rows.Item(1).Insert()
$legend.Range("A1").Value = "Column1"
$legend.Range("B1").Value = "Column2"

# Turn A1:B7 into an Excel Table, using the row we just wrote as headers.
$table = $legend.ListObjects.Add(1, $legend.Range("A1:B7"), $null, 1)
$table.Name = "Table1"

[void]$legend.Range("A1:B7").Select()

# Restore the sheet that was active before the edit.
$dataSheet.Activate()
